$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" date column (C) for rows 2-23 from 2023-09-08 (45177) to 2023-09-09 (45178)
for ($row = 2; $row -le 23; $row++) {
    $ws.Cells.Item($row, 3).Value = 45178
}
